# Maj dico de données : ajout de la classe répondre
#
# Adds a new data row (row 30) to the "Feuille 1" dictionary sheet,
# describing the new "Répondre" entity/class, reusing the same visual
# layout (styles + merges) as the existing rows, and updates the
# selection/view state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30 content -----------------------------------------------------
# Row 24 carries the exact same visual pattern we need for row 30
# (entity-name cell spanning A:B, then data/meaning/type/constraint
# columns). Merge first so the subsequent format-only paste lands the
# merge-aware border set (matches the existing rows' cellXfs exactly),
# then paste the formatting, then fill in the values.

$ws.Range("A30:B30").Merge()
$ws.Range("C30:D30").Merge()
$ws.Range("E30:F30").Merge()

$ws.Range("A24:I24").Copy()
$ws.Range("A30:I30").PasteSpecial(-4122)

$ws.Range("A30").Value = "Répondre"
$ws.Range("C30").Value = "idComm"
$ws.Range("E30").Value = "identifiant commentaire sujet"
$ws.Range("G30").Value = "AUTO_INCREMENT"
$ws.Range("I30").Value = "clé primaire"

# --- View / selection state ---------------------------------------------
$ws.Range("H30").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
